# Split the Title, Author and Abstract paragraphs' single run of text into
# one run per word plus one run per inter-word space, leaving the text and
# paragraph formatting unchanged.

$d = $word.ActiveDocument

function Escape-Xml($s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

function Build-RunsXml($text) {
    # Split on single spaces, keeping each space as its own "word" too, so
    # "Tom Coleman" -> "Tom", " ", "Coleman".
    $parts = $text -split ' '
    $inner = ""
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($i -gt 0) {
            $inner += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
        }
        $word = $parts[$i]
        if ($word -ne "") {
            $esc = Escape-Xml $word
            $inner += '<w:r><w:t xml:space="preserve">' + $esc + '</w:t></w:r>'
        }
    }
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $inner + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Split-ParagraphIntoWordRuns($paragraph) {
    $text = $paragraph.Range.Text
    # Range.Text includes the trailing paragraph-mark character(s) (CR and/or
    # cell/row marks); strip any trailing control characters so we only
    # split the visible text into words.
    $text = $text.TrimEnd([char]13, [char]7)
    $start = $paragraph.Range.Start
    $end = $paragraph.Range.End - 1   # exclude the trailing paragraph mark
    $target = $d.Range($start, $end)
    $xml = Build-RunsXml $text
    $target.InsertXML($xml)
}

Split-ParagraphIntoWordRuns $d.Paragraphs.Item(1)   # Title: "Questions: Using the quadratic formula"
Split-ParagraphIntoWordRuns $d.Paragraphs.Item(2)   # Author: "Tom Coleman"
Split-ParagraphIntoWordRuns $d.Paragraphs.Item(4)   # Abstract: "A selection of questions on using the quadratic formula."
